# summer 24 week 1 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.31
$ws.Range("B3").Value = 1.53
$ws.Range("E3").Value = 1.33
$ws.Range("F3").Value = 1.2
$ws.Range("F4").Value = 1.12
$ws.Range("C5").Value = 1.34
$ws.Range("C6").Value = 1.48
$ws.Range("D6").Value = 1.48
$ws.Range("G6").Value = 0.98
$ws.Range("F7").Value = 1.48
